# Auto-generated Excel COM-interop script
# Applies updated market-data values (H..N columns) across multiple sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 891.6667
$ws.Range("I6").Value = 317.57144
$ws.Range("J6").Value = 2901
$ws.Range("K6").Value = 952.71432
$ws.Range("L6").Value = 8703
$ws.Range("M6").Value = -840.71432
$ws.Range("N6").Value = -8927

$ws.Range("H8").Value = 38.166668
$ws.Range("I8").Value = 38.166668
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 114.500004
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 24.499996
$ws.Range("N8").ClearContents()

$ws.Range("H38").Value = 441.14285
$ws.Range("I38").Value = 217.6
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 652.8
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -280.8
$ws.Range("N38").Value = -3744

$ws.Range("H74").Value = 3770
$ws.Range("I74").Value = 3481.25
$ws.Range("J74").Value = 3935
$ws.Range("K74").Value = 3481.25
$ws.Range("L74").Value = 3935
$ws.Range("M74").Value = -2545.25
$ws.Range("N74").Value = -5807

$ws.Range("H77").Value = 3770
$ws.Range("I77").Value = 3481.25
$ws.Range("J77").Value = 3935
$ws.Range("K77").Value = 17406.25
$ws.Range("L77").Value = 19675
$ws.Range("M77").Value = -12726.25
$ws.Range("N77").Value = -29035

$ws.Range("H137").Value = 1088.5254
$ws.Range("I137").Value = 494.7619
$ws.Range("J137").Value = 2555.4707
$ws.Range("K137").Value = 1484.2857
$ws.Range("L137").Value = 7666.4121
$ws.Range("M137").Value = 1065.7143
$ws.Range("N137").Value = -12766.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14792.12
$ws.Range("I32").Value = 14413.318
$ws.Range("J32").Value = 18622.223
$ws.Range("K32").Value = 14413.318
$ws.Range("L32").Value = 18622.223
$ws.Range("M32").Value = -14126.318
$ws.Range("N32").Value = -19196.223

$ws.Range("H45").Value = 30304236
$ws.Range("I45").Value = 37038188
$ws.Range("J45").Value = 1450
$ws.Range("K45").Value = 37038188
$ws.Range("L45").Value = 1450
$ws.Range("M45").Value = -37037811
$ws.Range("N45").Value = -2204

$ws.Range("H74").Value = 664.875
$ws.Range("I74").Value = 468.7857
$ws.Range("J74").Value = 939.4
$ws.Range("K74").Value = 468.7857
$ws.Range("L74").Value = 939.4
$ws.Range("M74").Value = 405.2143
$ws.Range("N74").Value = -2687.4

$ws.Range("H77").Value = 664.875
$ws.Range("I77").Value = 468.7857
$ws.Range("J77").Value = 939.4
$ws.Range("K77").Value = 2343.9285
$ws.Range("L77").Value = 4697
$ws.Range("M77").Value = 2024.0715
$ws.Range("N77").Value = -13433

$ws.Range("H132").Value = 3375.1128
$ws.Range("I132").Value = 3352.7886
$ws.Range("J132").Value = 3491.2
$ws.Range("K132").Value = 10058.3658
$ws.Range("L132").Value = 10473.6
$ws.Range("M132").Value = -7528.3658
$ws.Range("N132").Value = -15533.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6023.294
$ws.Range("I20").Value = 6749.75
$ws.Range("J20").Value = 4279.8
$ws.Range("K20").Value = 6749.75
$ws.Range("L20").Value = 4279.8
$ws.Range("M20").Value = -6502.75
$ws.Range("N20").Value = -4773.8

$ws.Range("H105").Value = 2175.3076
$ws.Range("I105").Value = 1989.7778
$ws.Range("J105").Value = 2273.5293
$ws.Range("K105").Value = 1989.7778
$ws.Range("L105").Value = 2273.5293
$ws.Range("M105").Value = -242.7778000000001
$ws.Range("N105").Value = -5767.5293

$ws.Range("H134").Value = 2583.353
$ws.Range("I134").Value = 2917.7856
$ws.Range("J134").Value = 2349.25
$ws.Range("K134").Value = 8753.356800000001
$ws.Range("L134").Value = 7047.75
$ws.Range("M134").Value = -6218.356800000001
$ws.Range("N134").Value = -12117.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3291891.2
$ws.Range("I132").Value = 1763.05
$ws.Range("J132").Value = 6947589.5
$ws.Range("K132").Value = 5289.15
$ws.Range("L132").Value = 20842768.5
$ws.Range("M132").Value = -2759.15
$ws.Range("N132").Value = -20847828.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 65.5
$ws.Range("I4").Value = 65.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 196.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -84.5
$ws.Range("N4").ClearContents()

$ws.Range("H10").Value = 95.25
$ws.Range("I10").Value = 95.25
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 285.75
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -146.75

$ws.Range("H15").Value = 721.3333
$ws.Range("I15").Value = 200
$ws.Range("J15").Value = 982
$ws.Range("K15").Value = 600
$ws.Range("L15").Value = 2946
$ws.Range("M15").Value = -460
$ws.Range("N15").Value = -3226

$ws.Range("H16").Value = 200
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -427

$ws.Range("H17").Value = 425
$ws.Range("I17").Value = 350
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 1050
$ws.Range("L17").Value = 1500
$ws.Range("M17").Value = -881
$ws.Range("N17").Value = -1838

$ws.Range("H23").Value = 64.92856999999999
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 72.416664
$ws.Range("K23").Value = 60
$ws.Range("L23").Value = 217.249992
$ws.Range("M23").Value = 175
$ws.Range("N23").Value = -687.249992

$ws.Range("H25").Value = 1100
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1100
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 3300
$ws.Range("M25").ClearContents()
$ws.Range("N25").Value = -3638

$ws.Range("H30").Value = 1100
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 1100
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 3300
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -3504

$ws.Range("H39").Value = 4000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 4000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588

$ws.Range("H40").Value = 65.375
$ws.Range("I40").Value = 69.14286
$ws.Range("J40").Value = 39
$ws.Range("K40").Value = 276.57144
$ws.Range("L40").Value = 156
$ws.Range("M40").Value = -207.57144
$ws.Range("N40").Value = -294

$ws.Range("H120").Value = 21895.227
$ws.Range("I120").Value = 11206
$ws.Range("J120").Value = 25039.117
$ws.Range("K120").Value = 33618
$ws.Range("L120").Value = 75117.351
$ws.Range("M120").Value = -28780
$ws.Range("N120").Value = -84793.351

$ws.Range("H131").Value = 1386095.6
$ws.Range("I131").Value = 12797.375
$ws.Range("J131").Value = 1536594.1
$ws.Range("K131").Value = 38392.125
$ws.Range("L131").Value = 4609782.300000001
$ws.Range("M131").Value = -33352.125
$ws.Range("N131").Value = -4619862.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 38721852
$ws.Range("I70").Value = 52072656
$ws.Range("J70").Value = 4520.8
$ws.Range("K70").Value = 52072656
$ws.Range("L70").Value = 4520.8
$ws.Range("M70").Value = -52072386
$ws.Range("N70").Value = -5060.8

$ws.Range("H73").Value = 38721852
$ws.Range("I73").Value = 52072656
$ws.Range("J73").Value = 4520.8
$ws.Range("K73").Value = 52072656
$ws.Range("L73").Value = 4520.8
$ws.Range("M73").Value = -52071720
$ws.Range("N73").Value = -6392.8

$ws.Range("H102").Value = 1257.0526
$ws.Range("I102").Value = 1092.1666
$ws.Range("J102").Value = 1539.7142
$ws.Range("K102").Value = 1092.1666
$ws.Range("L102").Value = 1539.7142
$ws.Range("M102").Value = 529.8334
$ws.Range("N102").Value = -4783.7142

$ws.Range("H113").Value = 27778796
$ws.Range("I113").Value = 83334020
$ws.Range("J113").Value = 1183.3334
$ws.Range("K113").Value = 83334020
$ws.Range("L113").Value = 1183.3334
$ws.Range("M113").Value = -83331850
$ws.Range("N113").Value = -5523.3334

$ws.Range("H126").Value = 953842.9
$ws.Range("I126").Value = 1681.0526
$ws.Range("J126").Value = 2084535
$ws.Range("K126").Value = 5043.1578
$ws.Range("L126").Value = 6253605
$ws.Range("M126").Value = -2573.1578
$ws.Range("N126").Value = -6258545

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 33335586
$ws.Range("I7").Value = 2216
$ws.Range("J7").Value = 50002270
$ws.Range("K7").Value = 2216
$ws.Range("L7").Value = 50002270
$ws.Range("M7").Value = -2104
$ws.Range("N7").Value = -50002494

$ws.Range("H122").Value = 2633.7307
$ws.Range("I122").Value = 2687.4666
$ws.Range("J122").Value = 2560.4546
$ws.Range("K122").Value = 8062.399800000001
$ws.Range("L122").Value = 7681.3638
$ws.Range("M122").Value = -5612.399800000001
$ws.Range("N122").Value = -12581.3638

$ws.Range("H126").Value = 33335586
$ws.Range("I126").Value = 2216
$ws.Range("J126").Value = 50002270
$ws.Range("K126").Value = 6648
$ws.Range("L126").Value = 150006810
$ws.Range("M126").Value = -4178
$ws.Range("N126").Value = -150011750

$ws.Range("H132").Value = 6462.927
$ws.Range("I132").Value = 10952.8
$ws.Range("J132").Value = 2186.8572
$ws.Range("K132").Value = 32858.39999999999
$ws.Range("L132").Value = 6560.571599999999
$ws.Range("M132").Value = -30328.39999999999
$ws.Range("N132").Value = -11620.5716

$ws.Range("H133").Value = 25866
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 25866
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 25866
$ws.Range("N133").Value = -30926

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1303.6227
$ws.Range("I132").Value = 1076.6818
$ws.Range("J132").Value = 2413.111
$ws.Range("K132").Value = 3230.0454
$ws.Range("L132").Value = 7239.333
$ws.Range("M132").Value = -700.0454
$ws.Range("N132").Value = -12299.333

$ws.Range("H136").Value = 1051.2538
$ws.Range("I136").Value = 1071.8422
$ws.Range("J136").Value = 1024.2759
$ws.Range("K136").Value = 3215.5266
$ws.Range("L136").Value = 3072.8277
$ws.Range("M136").Value = -665.5266000000001
$ws.Range("N136").Value = -8172.8277
Write-Output "Applied market-data updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
